$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10th column) to make room for "Ghi chú"
$ws.Columns.Item(10).Insert()

# Set header text for the new column
$ws.Range("J1").Value = "Ghi chú"

# Merge the new header cells like the other header columns (J1:J2)
$ws.Range("J1:J2").Merge()

# Adjust the new column width to match the target layout (same width as column A)
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# The new column's data row (row 3) should have no border, unlike the
# surrounding bordered table cells
$ws.Range("J3").Borders.LineStyle = -4142

# Update the sheet view: scroll so column E is the top-left visible column,
# and set the active selection to K4
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("K4").Select()
